$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Sheet: Summary (row 15 -> Pattern3-Data+News / llama-3.1-405b) ---
$ws1 = $wb.Worksheets.Item("Summary")
Set-TextCell $ws1.Range("D15") "¥1,000,689.80"
Set-TextCell $ws1.Range("E15") "¥+689.80"
Set-TextCell $ws1.Range("F15") "+0.07%"
Set-TextCell $ws1.Range("G15") "+1.35%"
$ws1.Range("H15").Value = -1.201
Set-TextCell $ws1.Range("J15") "58.3%"
Set-TextCell $ws1.Range("K15") "0.0058%"
Set-TextCell $ws1.Range("L15") "0.4422%"
$ws1.Range("M15").Value = 13
$ws1.Range("N15").Value = 13
Set-TextCell $ws1.Range("P15") "20260106"

# --- Sheet: Pattern3-Data+News (row 5 -> llama-3.1-405b) ---
$ws4 = $wb.Worksheets.Item("Pattern3-Data+News")
Set-TextCell $ws4.Range("D5") "¥1,000,689.80"
Set-TextCell $ws4.Range("E5") "¥+689.80"
Set-TextCell $ws4.Range("F5") "+0.07%"
Set-TextCell $ws4.Range("G5") "+1.35%"
$ws4.Range("H5").Value = -1.201
Set-TextCell $ws4.Range("J5") "58.3%"
Set-TextCell $ws4.Range("K5") "0.0058%"
Set-TextCell $ws4.Range("L5") "0.4422%"
$ws4.Range("M5").Value = 13
$ws4.Range("N5").Value = 13
Set-TextCell $ws4.Range("P5") "20260106"
